$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.168.19'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.784.84'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '''226.04'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''32.05'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("D10").Value = '''0.0687'
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").Value = '''0.0948'
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '2.042.64'
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").Value = '''10.95'
$ws.Range("E13").Value = '  -3.22%  '
$ws.Range("D14").Value = '1.773.53'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").Value = '''0.624'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '34.166.27'
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '''67.67'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").Value = '''245.64'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").Value = '0.0₃0794'
$ws.Range("E20").Value = '  +2.27%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").Value = '''161.54'
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = '''16.30'
$ws.Range("E27").Value = '  +0.15%  '
$ws.Range("D28").Value = '''0.114'
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("E29").Value = '  +0.27%  '
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").Value = '''3.74'
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("E33").Value = '  +3.44%  '
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").Value = '1.445.08'
$ws.Range("E35").Value = '  +2.93%  '
$ws.Range("E36").Value = '  +11.16%  '
$ws.Range("D37").Value = '''0.654'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("E38").Value = '  +1.14%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").Value = '''82.24'
$ws.Range("E40").Value = '  +2.73%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''14.17'
$ws.Range("E41").Value = '  +5.97%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '''2.38'
$ws.Range("E42").Value = '  +0.96%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.72'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''0.915'
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D45").Value = '''0.0517'
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").Value = '''1.08'
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("D48").Value = '1.942.27'
$ws.Range("E48").Value = '  -0.23%  '
$ws.Range("D49").Value = '''104.82'
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("E50").Value = '  -5.33%  '
$ws.Range("E51").Value = '  +0.21%  '
